$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1244.7693
$ws.Range("I137").Value = 1181.9546
$ws.Range("K137").Value = 3545.8638
$ws.Range("M137").Value = -995.8638000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 28497.676
$ws.Range("I2").Value = 1233.68
$ws.Range("J2").Value = 85297.664
$ws.Range("K2").Value = 1233.68
$ws.Range("L2").Value = 85297.664
$ws.Range("M2").Value = -1120.68
$ws.Range("N2").Value = -85523.664
$ws.Range("H32").Value = 27252.672
$ws.Range("I32").Value = 5327.0415
$ws.Range("J32").Value = 82643.734
$ws.Range("K32").Value = 5327.0415
$ws.Range("L32").Value = 82643.734
$ws.Range("M32").Value = -5040.0415
$ws.Range("N32").Value = -83217.734
$ws.Range("H74").Value = 2419.9211
$ws.Range("I74").Value = 1637.55
$ws.Range("J74").Value = 3289.2222
$ws.Range("K74").Value = 1637.55
$ws.Range("L74").Value = 3289.2222
$ws.Range("M74").Value = -763.55
$ws.Range("N74").Value = -5037.2222
$ws.Range("H77").Value = 2419.9211
$ws.Range("I77").Value = 1637.55
$ws.Range("J77").Value = 3289.2222
$ws.Range("K77").Value = 8187.75
$ws.Range("L77").Value = 16446.111
$ws.Range("M77").Value = -3819.75
$ws.Range("N77").Value = -25182.111
$ws.Range("H102").Value = 202836
$ws.Range("I102").Value = 502490
$ws.Range("K102").Value = 502490
$ws.Range("M102").Value = -500868
$ws.Range("H116").Value = 28497.676
$ws.Range("I116").Value = 1233.68
$ws.Range("J116").Value = 85297.664
$ws.Range("K116").Value = 1233.68
$ws.Range("L116").Value = 85297.664
$ws.Range("M116").Value = 1060.32
$ws.Range("N116").Value = -89885.664
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 28497.676
$ws.Range("I3").Value = 1233.68
$ws.Range("J3").Value = 85297.664
$ws.Range("K3").Value = 1233.68
$ws.Range("L3").Value = 85297.664
$ws.Range("M3").Value = -1119.68
$ws.Range("N3").Value = -85525.664
$ws.Range("H134").Value = 1769.8572
$ws.Range("I134").Value = 1727.1666
$ws.Range("K134").Value = 5181.4998
$ws.Range("M134").Value = -2646.4998
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14801.792
$ws.Range("I31").Value = 31377.031
$ws.Range("J31").Value = 2370.3635
$ws.Range("K31").Value = 31377.031
$ws.Range("L31").Value = 2370.3635
$ws.Range("M31").Value = -31082.031
$ws.Range("N31").Value = -2960.3635
$ws.Range("H34").Value = 14801.792
$ws.Range("I34").Value = 31377.031
$ws.Range("J34").Value = 2370.3635
$ws.Range("K34").Value = 31377.031
$ws.Range("L34").Value = 2370.3635
$ws.Range("M34").Value = -31175.031
$ws.Range("N34").Value = -2774.3635
$ws.Range("H58").Value = 14553.8
$ws.Range("I58").Value = 1725.6666
$ws.Range("J58").Value = 130007
$ws.Range("K58").Value = 1725.6666
$ws.Range("L58").Value = 130007
$ws.Range("M58").Value = -1522.6666
$ws.Range("N58").Value = -130413
$ws.Range("H134").Value = 1631.2222
$ws.Range("I134").Value = 1235.8462
$ws.Range("K134").Value = 3707.5386
$ws.Range("M134").Value = -1172.5386
$ws.Range("H136").Value = 14553.8
$ws.Range("I136").Value = 1725.6666
$ws.Range("J136").Value = 130007
$ws.Range("K136").Value = 5176.9998
$ws.Range("L136").Value = 390021
$ws.Range("M136").Value = -2626.9998
$ws.Range("N136").Value = -395121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 18484.166
$ws.Range("I80").Value = 951
$ws.Range("J80").Value = 27250.75
$ws.Range("K80").Value = 2853
$ws.Range("L80").Value = 81752.25
$ws.Range("M80").Value = -1917
$ws.Range("N80").Value = -83624.25
$ws.Range("H83").Value = 18484.166
$ws.Range("I83").Value = 951
$ws.Range("J83").Value = 27250.75
$ws.Range("K83").Value = 8559
$ws.Range("L83").Value = 245256.75
$ws.Range("M83").Value = -3879
$ws.Range("N83").Value = -254616.75
$ws.Range("H129").Value = 1958.6666
$ws.Range("I129").Value = 580
$ws.Range("J129").Value = 2488.923
$ws.Range("K129").Value = 1740
$ws.Range("L129").Value = 7466.768999999999
$ws.Range("M129").Value = 3260
$ws.Range("N129").Value = -17466.769
$ws.Range("H130").Value = 1039.8
$ws.Range("I130").Value = 751.3333
$ws.Range("J130").Value = 1472.5
$ws.Range("K130").Value = 2253.9999
$ws.Range("L130").Value = 4417.5
$ws.Range("M130").Value = 2766.0001
$ws.Range("N130").Value = -14457.5
$ws.Range("H131").Value = 1385.3818
$ws.Range("J131").Value = 1399.102
$ws.Range("L131").Value = 4197.306
$ws.Range("N131").Value = -14277.306
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5185.696
$ws.Range("I122").Value = 5342.6924
$ws.Range("K122").Value = 16028.0772
$ws.Range("M122").Value = -13578.0772
$ws.Range("H126").Value = 2199.3333
$ws.Range("I126").Value = 2598.8
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 7796.400000000001
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -5326.400000000001
$ws.Range("N126").Value = -10040
$ws.Range("H132").Value = 3303.125
$ws.Range("I132").Value = 2693.7
$ws.Range("K132").Value = 8081.099999999999
$ws.Range("M132").Value = -5551.099999999999
$ws.Range("H136").Value = 47426.11
$ws.Range("J136").Value = 47426.11
$ws.Range("L136").Value = 142278.33
$ws.Range("N136").Value = -147378.33
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 127037.5
$ws.Range("I16").Value = 201060
$ws.Range("J16").Value = 3666.6667
$ws.Range("K16").Value = 201060
$ws.Range("L16").Value = 3666.6667
$ws.Range("M16").Value = -200890
$ws.Range("N16").Value = -4006.6667
$ws.Range("H68").Value = 2455.9
$ws.Range("I68").Value = 1455
$ws.Range("J68").Value = 4314.7144
$ws.Range("K68").Value = 1455
$ws.Range("L68").Value = 4314.7144
$ws.Range("M68").Value = -706
$ws.Range("N68").Value = -5812.7144
$ws.Range("H71").Value = 2455.9
$ws.Range("I71").Value = 1455
$ws.Range("J71").Value = 4314.7144
$ws.Range("K71").Value = 7275
$ws.Range("L71").Value = 21573.572
$ws.Range("M71").Value = -3531
$ws.Range("N71").Value = -29061.572
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 822.8125
$ws.Range("I136").Value = 496.52173
$ws.Range("J136").Value = 1656.6666
$ws.Range("K136").Value = 1489.56519
$ws.Range("L136").Value = 4969.9998
$ws.Range("M136").Value = 1060.43481
$ws.Range("N136").Value = -10069.9998

Write-Host "Applied all updates"